$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the cryptos.xlsx price/volume refresh diff.
# D-column values that look like plain decimals (e.g. "586.46") are written
# with a leading apostrophe so Excel stores them as text (matching the
# original inlineStr cells) instead of silently parsing them as numbers.
# D-column values containing two dots (e.g. "66.988.15") are never
# number-like, so they can be assigned directly.

$ws.Range("D2").Value = '66.988.15'
$ws.Range("E2").Value = '  -0.63%  '
$ws.Range("D3").Value = '3.522.81'
$ws.Range("E3").Value = '  +0.87%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''586.46'
$ws.Range("E5").Value = '  -1.74%  '
$ws.Range("D6").Value = '''178.09'
$ws.Range("E6").Value = '  -0.73%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("D8").Value = '''0.603'
$ws.Range("E8").Value = '  -0.06%  '
$ws.Range("D9").Value = '3.521.71'
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("E10").Value = '  -1.79%  '
$ws.Range("E11").Value = '  -1.68%  '
$ws.Range("E12").Value = '  -2.47%  '
$ws.Range("D13").Value = '4.135.82'
$ws.Range("E13").Value = '  +0.92%  '
$ws.Range("D14").Value = '''30.64'
$ws.Range("E14").Value = '  -5.35%  '
$ws.Range("D16").Value = '66.963.03'
$ws.Range("E16").Value = '  -0.69%  '
$ws.Range("E17").Value = '  -1.99%  '
$ws.Range("D18").Value = '3.523.54'
$ws.Range("E18").Value = '  +1.02%  '
$ws.Range("E19").Value = '  -3.16%  '
$ws.Range("D20").Value = '''14.08'
$ws.Range("E20").Value = '  -1.59%  '
$ws.Range("D21").Value = '''383.68'
$ws.Range("E21").Value = '  -1.52%  '
$ws.Range("E22").Value = '  -1.11%  '
$ws.Range("D23").Value = '''0.552'
$ws.Range("E23").Value = '  +2.00%  '
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("D25").Value = '''72.54'
$ws.Range("E25").Value = '  -2.11%  '
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("E27").Value = '  -0.01%  '
$ws.Range("E28").Value = '  -4.59%  '
$ws.Range("E29").Value = '  -1.46%  '
$ws.Range("D30").Value = '''0.999'
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("D31").Value = '''24.75'
$ws.Range("D32").Value = '''5.93'
$ws.Range("E32").Value = '  -4.85%  '
$ws.Range("E33").Value = '  -1.98%  '
$ws.Range("E34").Value = '  -5.19%  '
$ws.Range("D35").Value = '''7.30'
$ws.Range("E35").Value = '  -1.32%  '
$ws.Range("D36").Value = '''0.999'
$ws.Range("E36").Value = '  -0.06%  '
$ws.Range("E37").Value = '  -0.90%  '
$ws.Range("D38").Value = '''30.04'
$ws.Range("E38").Value = '  +14.30%  '
$ws.Range("D39").Value = '''161.34'
$ws.Range("E39").Value = '  -2.24%  '
$ws.Range("D40").Value = '''0.898'
$ws.Range("E40").Value = '  +3.21%  '
$ws.Range("E41").Value = '  -4.14%  '
$ws.Range("D42").Value = '''6.63'
$ws.Range("E42").Value = '  -2.56%  '
$ws.Range("E43").Value = '  -2.37%  '
$ws.Range("D44").Value = '''2.58'
$ws.Range("E44").Value = '  -8.04%  '
$ws.Range("D45").Value = '2.741.64'
$ws.Range("E45").Value = '  -3.61%  '
$ws.Range("E46").Value = '  -2.25%  '
$ws.Range("D47").Value = '''25.34'
$ws.Range("E47").Value = '  -6.29%  '
$ws.Range("E48").Value = '  -2.17%  '
$ws.Range("E49").Value = '  -0.70%  '
$ws.Range("D50").Value = '''324.19'
$ws.Range("E50").Value = '  -3.52%  '
$ws.Range("E51").Value = '  -3.08%  '
